$wb = $excel.ActiveWorkbook

# --- StudentsMapping sheet: organization name assignments ---
$ws2 = $wb.Worksheets.Item("StudentsMapping")
$ws2.Range("C2").Value = "Providence Health Network"
$ws2.Range("D2").Value = "Cedars-Sinai - Neurosciences"
$ws2.Range("E2").Value = "Keck VIO - COBI"
$ws2.Range("F2").Value = "Optum CF - Patient XP"
$ws2.Range("D4").Value = "Rancho Los Amigos NRC"
$ws2.Range("E4").Value = "CHLA - Anesthesia&CCM"
$ws2.Range("F4").Value = "Providence Health Network"
$ws2.Range("D5").Value = "Emanate Health"
$ws2.Range("E5").Value = "Cedars-Sinai - Neurosciences"
$ws2.Range("E6").Value = "Verdugo Hills Hospital"
$ws2.Range("F6").Value = "Keck VIO - COBI"
$ws2.Range("E7").Value = "Providence Health Network"
$ws2.Range("D8").Value = "CHLA - Anesthesia&CCM"
$ws2.Range("E8").Value = "City of Hope - CMO"
$ws2.Range("F8").Value = "Rancho Los Amigos NRC"
$ws2.Range("C9").Value = "Cedars-Sinai - Neurosciences"
$ws2.Range("D9").Value = "Keck VIO - COBI"
$ws2.Range("E9").Value = "Optum CF - Patient XP"
$ws2.Range("F9").Value = "St.Johns-PhysPartners "
$ws2.Range("F12").Value = "Torrance Memorial"

# --- OrganizationMapping sheet: student name assignments ---
$ws3 = $wb.Worksheets.Item("OrganizationMapping")
$ws3.Range("C2").Value = "Raashi Subramanya"
$ws3.Range("D2").Value = "Oceana Hanner"
$ws3.Range("E2").Value = "Fahima Gohil"
$ws3.Range("D3").Value = "Stanley Ibe"
$ws3.Range("E3").Value = "Daniela Ahumada"
$ws3.Range("D4").Value = ""
$ws3.Range("E4").Value = "Stanley Ibe"
$ws3.Range("D5").Value = "Fahima Gohil"
$ws3.Range("E5").Value = ""
$ws3.Range("F6").Value = ""
$ws3.Range("D8").Value = "Raashi Subramanya"
$ws3.Range("E8").Value = "Oceana Hanner"
$ws3.Range("F8").Value = "Julia Orozco"
$ws3.Range("E10").Value = "Raashi Subramanya"
$ws3.Range("F10").Value = "Oceana Hanner"
$ws3.Range("C11").Value = "Oceana Hanner"
$ws3.Range("E11").Value = "Emma Crusinberry"
$ws3.Range("F11").Value = "Daniela Ahumada"
$ws3.Range("D12").Value = "Daniela Ahumada"
$ws3.Range("F12").Value = "Stanley Ibe"
$ws3.Range("D13").Value = ""
$ws3.Range("E13").Value = ""
$ws3.Range("F14").Value = "Raashi Subramanya"
$ws3.Range("F15").Value = "Russelle Chang"
$ws3.Range("E16").Value = "Julia Orozco"
